$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Enable (Runmode = Y) the 3 newly added search test cases (rows 43-45)
$ws.Range("D43").Value = "Y"
$ws.Range("D44").Value = "Y"
$ws.Range("D45").Value = "Y"

# Update Results column to reflect the latest run outcome
$ws.Range("E23").Value = "FAIL"
$ws.Range("E24").Value = "FAIL"
$ws.Range("E25").Value = "FAIL"
$ws.Range("E30").Value = "FAIL"
$ws.Range("E43").Value = "FAIL"
$ws.Range("E44").Value = "PASS"
$ws.Range("E45").Value = "PASS"

# Reflect the final cursor/view position (last edited cell)
$ws.Activate()
$ws.Range("D44").Select()
